$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value2 = "43.749.88"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value2 = "  +1.12%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value2 = "2.246.14"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value2 = "  +0.49%  "
$ws.Cells.Item(4, 5).Value2 = "  +0.08%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value2 = "323.33"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value2 = "  +2.57%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value2 = "101.83"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value2 = "  +0.02%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.579"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value2 = "  -1.37%  "
$ws.Cells.Item(8, 5).Value2 = "  +0.04%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.556"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value2 = "  -0.92%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value2 = "37.40"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value2 = "  +0.89%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.0829"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value2 = "  +0.43%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value2 = "7.72"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value2 = "  +0.51%  "
$ws.Cells.Item(13, 5).Value2 = "  -1.92%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value2 = "2.588.04"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value2 = "  +0.44%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.859"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value2 = "  -0.21%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value2 = "14.22"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value2 = "  -0.95%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value2 = "2.247.46"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value2 = "  +0.36%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value2 = "43.686.97"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value2 = "  +1.17%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value2 = "13.82"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value2 = "  -3.53%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.0₃0987"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value2 = "  +2.62%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value2 = "6.57"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value2 = "  +1.20%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value2 = "65.28"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value2 = "  -0.50%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.18"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value2 = "  +0.19%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value2 = "236.41"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value2 = "  -1.02%  "
$ws.Cells.Item(25, 5).Value2 = "  +1.83%  "
$ws.Cells.Item(26, 5).Value2 = "  +0.12%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value2 = "10.13"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value2 = "  +1.36%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value2 = "2.20"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value2 = "  -1.62%  "
$ws.Cells.Item(29, 5).Value2 = "  +7.12%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value2 = "6.31"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value2 = "  -1.10%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value2 = "160.62"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value2 = "  +3.93%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value2 = "20.19"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value2 = "  -1.69%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.0854"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value2 = "  -2.21%  "
$ws.Cells.Item(34, 5).Value2 = "  -2.76%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.18"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value2 = "  -0.19%  "
$ws.Cells.Item(36, 5).Value2 = "  +9.14%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value2 = "1.95"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value2 = "  +0.63%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.120"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value2 = "  -1.83%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.79"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value2 = "  +1.79%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value2 = "4.29"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value2 = "  -3.77%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value2 = "15.75"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value2 = "  +21.96%  "
$ws.Cells.Item(42, 5).Value2 = "  -1.61%  "
$ws.Cells.Item(43, 5).Value2 = "  +0.28%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value2 = "1.820.30"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value2 = "  +1.10%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.201"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value2 = "  -2.38%  "
$ws.Cells.Item(46, 2).Value2 = "BitcoinSV"
$ws.Cells.Item(46, 3).Value2 = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value2 = "82.88"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value2 = "  -7.02%  "
$ws.Cells.Item(47, 2).Value2 = "Stacks"
$ws.Cells.Item(47, 3).Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value2 = "1.73"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value2 = "  +6.63%  "
$ws.Cells.Item(48, 5).Value2 = "  -1.98%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value2 = "74.64"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value2 = "  -3.16%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value2 = "58.89"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value2 = "  -0.57%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value2 = "103.67"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value2 = "  +0.24%  "
